$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '2026-02-08 20:18:38'
$ws.Range("I2").Value = '5.4 mm'
$ws.Range("E3").Value = '2026-02-08 20:18:40'
$ws.Range("E4").Value = '2026-02-08 20:18:43'
$ws.Range("J4").Value = '1002.2 hPa'
$ws.Range("O4").Value = '10.3 °C'
$ws.Range("E5").Value = '2026-02-08 20:18:45'
$ws.Range("G5").Value = '120 cm'
$ws.Range("I5").Value = '5.3 mm'
$ws.Range("E6").Value = '2026-02-08 20:18:48'
$ws.Range("J6").Value = '1002.0 hPa'
$ws.Range("E7").Value = '2026-02-08 20:18:50'
$ws.Range("J7").Value = '1002.4 hPa'
$ws.Range("E8").Value = '2026-02-08 20:18:53'
$ws.Range("J8").Value = '1002.3 hPa'
$ws.Range("E9").Value = '2026-02-08 20:18:55'
$ws.Range("H9").Value = '''69%'
$ws.Range("E10").Value = '2026-02-08 20:18:58'
$ws.Range("E11").Value = '2026-02-08 20:19:00'
$ws.Range("E12").Value = '2026-02-08 20:19:02'
$ws.Range("E13").Value = '2026-02-08 20:19:05'
$ws.Range("H13").Value = '''78%'
$ws.Range("J13").Value = '1003.7 hPa'
$ws.Range("E14").Value = '2026-02-08 20:19:07'
$ws.Range("E15").Value = '2026-02-08 20:19:10'
$ws.Range("O15").Value = '9.9 °C'
$ws.Range("E16").Value = '2026-02-08 20:19:12'
$ws.Range("H16").Value = '''81%'
$ws.Range("I16").Value = '3.5 mm'
$ws.Range("E17").Value = '2026-02-08 20:19:14'
$ws.Range("H17").Value = '''100%'
$ws.Range("E18").Value = '2026-02-08 20:19:17'
$ws.Range("I18").Value = '0.5 mm'
$ws.Range("J18").Value = '1002.4 hPa'
$ws.Range("E19").Value = '2026-02-08 20:19:19'
$ws.Range("E20").Value = '2026-02-08 20:19:22'
$ws.Range("I20").Value = '8.8 mm'
$ws.Range("E21").Value = '2026-02-08 20:19:24'
$ws.Range("H21").Value = '''79%'
$ws.Range("J21").Value = '1003.1 hPa'
$ws.Range("E22").Value = '2026-02-08 20:19:26'
$ws.Range("E23").Value = '2026-02-08 20:19:29'
$ws.Range("I23").Value = '5.2 mm'
$ws.Range("E24").Value = '2026-02-08 20:19:31'
$ws.Range("J24").Value = '1003.7 hPa'
$ws.Range("E25").Value = '2026-02-08 20:19:34'
$ws.Range("H25").Value = '''79%'
$ws.Range("O25").Value = '-3.0 °C'
$ws.Range("E26").Value = '2026-02-08 20:19:36'
$ws.Range("H26").Value = '''69%'
$ws.Range("J26").Value = '1001.4 hPa'
$ws.Range("E27").Value = '2026-02-08 20:19:39'
$ws.Range("H27").Value = '''88%'
$ws.Range("E28").Value = '2026-02-08 20:19:41'
$ws.Range("J28").Value = '1002.0 hPa'
$ws.Range("E29").Value = '2026-02-08 20:19:44'
$ws.Range("E30").Value = '2026-02-08 20:19:46'
$ws.Range("J30").Value = '1002.4 hPa'
$ws.Range("E31").Value = '2026-02-08 20:19:48'
$ws.Range("I31").Value = '0.6 mm'
$ws.Range("J31").Value = '1001.5 hPa'
$ws.Range("N31").Value = '7.7 °C 19:59 TU'
$ws.Range("O31").Value = '9.7 °C'
$ws.Range("E32").Value = '2026-02-08 20:19:51'
$ws.Range("E33").Value = '2026-02-08 20:19:53'
$ws.Range("O33").Value = '2.9 °C'
$ws.Range("E34").Value = '2026-02-08 20:19:56'
$ws.Range("H34").Value = '''73%'
$ws.Range("E35").Value = '2026-02-08 20:19:58'
$ws.Range("J35").Value = '1004.5 hPa'
$ws.Range("E36").Value = '2026-02-08 20:20:00'
$ws.Range("J36").Value = '1002.5 hPa'
$ws.Range("E37").Value = '2026-02-08 20:20:03'
$ws.Range("J37").Value = '1003.3 hPa'
$ws.Range("E38").Value = '2026-02-08 20:20:05'
$ws.Range("H38").Value = '''76%'
$ws.Range("I38").Value = '4.6 mm'
$ws.Range("O38").Value = '9.5 °C'
$ws.Range("E39").Value = '2026-02-08 20:20:08'
$ws.Range("E40").Value = '2026-02-08 20:20:10'
$ws.Range("J40").Value = '1003.7 hPa'
$ws.Range("O40").Value = '5.8 °C'
$ws.Range("E41").Value = '2026-02-08 20:20:13'
$ws.Range("H41").Value = '''69%'
$ws.Range("J41").Value = '1002.5 hPa'
$ws.Range("O41").Value = '12.1 °C'
$ws.Range("E42").Value = '2026-02-08 20:20:15'
$ws.Range("E43").Value = '2026-02-08 20:20:17'
$ws.Range("H43").Value = '''85%'
$ws.Range("E44").Value = '2026-02-08 20:20:20'
$ws.Range("E45").Value = '2026-02-08 20:20:22'
$ws.Range("J45").Value = '1004.6 hPa'
$ws.Range("E46").Value = '2026-02-08 20:20:25'
$ws.Range("J46").Value = '1004.1 hPa'
$ws.Range("O46").Value = '9.7 °C'
